# Add 4 new sheets, one per teammate, each pre-populated with the same
# header row as the original schedule sheet, and rename the original
# sheet to reflect its owner as well.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Copy the header row (A1:F1) from the original sheet so the new sheets
# share the same styling / shared-string values.
$ws1.Range("A1:F1").Copy()

$names = @(
    "Joseph Pak (Poducer and Design)",
    "Wyatt(Artist)",
    "Emanuel(designer)",
    "Raphael Brown(Programmer)"
)

$prev = $ws1
foreach ($n in $names) {
    $newSheet = $wb.Worksheets.Add($null, $prev)
    $newSheet.Name = $n
    $newSheet.Range("A1").PasteSpecial()
    $prev = $newSheet
}

# Rename the original sheet to identify its owner, matching the others.
$ws1.Name = "Marcus Queiro(Sound Effects)"

# Make the second tab (Joseph Pak's sheet) the active one, as in the
# final workbook.
$wb.Worksheets.Item(2).Activate()
